$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.009.96"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "3.420.92"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.19"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.01"
$ws.Range("E6").Value = "  +5.28%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.02"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("D12").Value = "4.004.54"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.35"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "3.416.96"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "62.042.79"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.45"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.96"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.14"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "3.561.90"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.56"
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.95"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "31.03"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "168.14"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "3.457.21"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.72"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.781"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "2.558.64"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.09"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("E51").Value = "  +0.04%  "
